# Updated all state election data
# Fills in the newly-modelled percentile columns (C:G) for the rows that
# previously only had an "Actual" value (column B) on four sheets, and
# updates the saved selection / active-sheet state to match.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. New percentile data (columns C:G, rows 3-8) for the affected sheets.
# -----------------------------------------------------------------------

$sheetData = @{
    "2013 Fed" = @{
        3 = @{ 3 = 44.06;                 4 = 44.6;                  5 = 45.22;                 6 = 45.85;                 7 = 46.35 }
        4 = @{ 3 = 32.57;                 4 = 33.119999999999997;    5 = 33.729999999999997;    6 = 34.36;                 7 = 34.83 }
        5 = @{ 3 = 8.86;                  4 = 9.3800000000000008;    5 = 10;                    6 = 10.63;                 7 = 11.15 }
        6 = @{ 3 = 11.27;                 4 = 11.8;                  5 = 12.44;                 6 = 13.03;                 7 = 13.54 }
        7 = @{ 3 = 3.3;                   4 = 3.96;                  5 = 4.76;                  6 = 5.58;                  7 = 6.24 }
    }
    "2010 Fed" = @{
        3 = @{ 3 = 41.1;                  4 = 41.68;                 5 = 42.42;                 6 = 43.13;                 7 = 43.7 }
        4 = @{ 3 = 37.03;                 4 = 37.61;                 5 = 38.32;                 6 = 39.04;                 7 = 39.700000000000003 }
        5 = @{ 3 = 11.63;                 4 = 12.26;                 5 = 12.99;                 6 = 13.72;                 7 = 14.3 }
        6 = @{ 3 = 4.6500000000000004;    4 = 5.27;                  5 = 6;                     6 = 6.7;                   7 = 7.29 }
    }
    "2007 Fed" = @{
        3 = @{ 3 = 40.11;                 4 = 40.69;                 5 = 41.44;                 6 = 42.2;                  7 = 42.79 }
        4 = @{ 3 = 44.06;                 4 = 44.69;                 5 = 45.43;                 6 = 46.18;                 7 = 46.77 }
        5 = @{ 3 = 6.45;                  4 = 7.14;                  5 = 7.99;                  6 = 8.83;                  7 = 9.5299999999999994 }
        6 = @{ 3 = 3.69;                  4 = 4.3899999999999997;    5 = 5.26;                  6 = 6.15;                  7 = 6.9 }
    }
    "2019 NSW" = @{
        3 = @{ 3 = 38.33;                 4 = 39.26;                 5 = 40.39;                 6 = 41.49;                 7 = 42.39 }
        4 = @{ 3 = 32.4;                  4 = 33.28;                 5 = 34.43;                 6 = 35.590000000000003;    7 = 36.49 }
        5 = @{ 3 = 8.65;                  4 = 9.51;                  5 = 10.63;                 6 = 11.75;                 7 = 12.7 }
        6 = @{ 3 = 13.35;                 4 = 14.28;                 5 = 15.37;                 6 = 16.489999999999998;    7 = 17.41 }
        7 = @{ 3 = 1.1200000000000001;    4 = 2.12;                  5 = 4.03;                  6 = 5.64;                  7 = 7.05 }
        8 = @{ 3 = 2.89;                  4 = 3.91;                  5 = 5.22;                  6 = 6.5;                   7 = 7.55 }
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    foreach ($r in $rows.Keys) {
        $cols = $rows[$r]
        foreach ($c in $cols.Keys) {
            $ws.Cells.Item($r, $c).Value = $cols[$c]
        }
    }
}

# -----------------------------------------------------------------------
# 2. Saved selection on each touched sheet (also drives tabSelected / the
#    workbook's active-tab bookkeeping as each sheet is selected in turn).
# -----------------------------------------------------------------------

[void]$wb.Worksheets.Item("2019 Fed").Range("D10").Select()
[void]$wb.Worksheets.Item("2013 Fed").Range("E5").Select()
[void]$wb.Worksheets.Item("2010 Fed").Range("E9").Select()
[void]$wb.Worksheets.Item("2007 Fed").Range("G9").Select()

# "2019 NSW" is selected last so it ends up as the active / tabSelected sheet.
$ws6 = $wb.Worksheets.Item("2019 NSW")
[void]$ws6.Activate()
[void]$ws6.Range("C13").Select()

# -----------------------------------------------------------------------
# 3. Scroll the workbook tab strip so sheet index 3 ("2010 Fed") is the
#    first visible tab, matching the saved window scroll position.
# -----------------------------------------------------------------------

$window = $excel.ActiveWindow
$window.ScrollWorkbookTabs(1, 3) | Out-Null
